$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells whose new text would otherwise be auto-coerced
# to a number by Excel (losing formatting like trailing zeros), so force
# the cell format to Text ("@") before assigning the string value.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D11",
    "D12",
    "D13",
    "D16",
    "D18",
    "D21",
    "D22",
    "D24",
    "D26",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.928.65"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "3.410.22"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "408.81"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "128.76"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +5.87%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.730"
$ws.Range("E9").Value = "  +5.24%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").Value = "42.65"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "9.14"
$ws.Range("E12").Value = "  +8.80%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  +36.80%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "3.953.07"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "21.15"
$ws.Range("E16").Value = "  +6.67%  "
$ws.Range("D17").Value = "3.402.99"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "12.43"
$ws.Range("E18").Value = "  +7.81%  "
$ws.Range("E19").Value = "  +6.78%  "
$ws.Range("D20").Value = "61.898.21"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "446.65"
$ws.Range("E21").Value = "  +43.03%  "
$ws.Range("D22").Value = "91.20"
$ws.Range("E22").Value = "  +8.08%  "
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "13.10"
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("E25").Value = "  +3.43%  "
$ws.Range("D26").Value = "9.28"
$ws.Range("E26").Value = "  +14.73%  "
$ws.Range("D27").Value = "32.91"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "7.61"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "12.05"
$ws.Range("E30").Value = "  +5.91%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "0.170"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "42.62"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "0.0501"
$ws.Range("E36").Value = "  +3.77%  "
$ws.Range("D37").Value = "53.81"
$ws.Range("E37").Value = "  +3.92%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E40").Value = "  +6.94%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("D42").Value = "0.317"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").Value = "142.60"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "4.23"
$ws.Range("E44").Value = "  +7.96%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.55"
$ws.Range("E45").Value = "  +15.08%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "1.99"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "22.29"
$ws.Range("E48").Value = "  +4.70%  "
$ws.Range("D49").Value = "0.145"
$ws.Range("E49").Value = "  +21.02%  "
$ws.Range("D50").Value = "2.13"
$ws.Range("E50").Value = "  +8.91%  "
$ws.Range("D51").Value = "3.757.53"
$ws.Range("E51").Value = "  -0.61%  "
